# Replace the small Name/Address/Position roster with a single-column
# "Serial" list of WDPE059A## tags (WDPE059A54 .. WDPE059A80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting prep -------------------------------------------------
# Row 3 (A3) already carries the plain bordered/centered style that every
# data row in the new layout should use. Clone that formatting down onto
# A2 and A4:A28 *before* we touch values, so the engine reuses the
# existing style slot instead of minting new ones.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A4:A28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Columns B and C are no longer used at all - wipe them completely
# (content + formatting) so no stray styled-but-empty cells remain.
$ws.Range("B1:C3").Clear()

# --- New data ----------------------------------------------------------
$ws.Range("A1").Value = "Serial"

$serials = @(
    "WDPE059A54","WDPE059A55","WDPE059A56","WDPE059A57","WDPE059A58",
    "WDPE059A59","WDPE059A60","WDPE059A61","WDPE059A62","WDPE059A63",
    "WDPE059A64","WDPE059A65","WDPE059A66","WDPE059A67","WDPE059A68",
    "WDPE059A69","WDPE059A70","WDPE059A71","WDPE059A72","WDPE059A73",
    "WDPE059A74","WDPE059A75","WDPE059A76","WDPE059A77","WDPE059A78",
    "WDPE059A79","WDPE059A80"
)

for ($i = 0; $i -lt $serials.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $serials[$i]
}

# Matches the author's final on-screen selection when the file was saved.
$ws.Range("C36").Select() | Out-Null
